$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so values like "1.000" or "16.00"
# are not auto-coerced into numbers and lose their exact textual representation.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.089.52'
$ws.Range("E2").Value = '  -3.86%  '

$ws.Range("D3").Value = '1.644.36'
$ws.Range("E3").Value = '  -3.66%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '307.52'
$ws.Range("E5").Value = '  -3.03%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").Value = '0.3903'
$ws.Range("E7").Value = '  -2.46%  '

$ws.Range("D8").Value = '0.3852'
$ws.Range("E8").Value = '  -4.80%  '

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("D10").Value = '1.348'
$ws.Range("E10").Value = '  -8.68%  '

$ws.Range("D11").Value = '49.02'
$ws.Range("E11").Value = '  -7.35%  '

$ws.Range("D12").Value = '0.08455'
$ws.Range("E12").Value = '  -4.14%  '

$ws.Range("D13").Value = '23.84'
$ws.Range("E13").Value = '  -8.63%  '

$ws.Range("D14").Value = '7.121'
$ws.Range("E14").Value = '  -4.99%  '

$ws.Range("D15").Value = '0.00001282'
$ws.Range("E15").Value = '  -5.31%  '

$ws.Range("E16").Value = '  -6.28%  '

$ws.Range("D17").Value = '1.641.96'
$ws.Range("E17").Value = '  -3.99%  '

$ws.Range("D18").Value = '94.85'
$ws.Range("E18").Value = '  -1.29%  '

$ws.Range("D19").Value = '0.06956'
$ws.Range("E19").Value = '  -3.36%  '

$ws.Range("D20").Value = '20.83'
$ws.Range("E20").Value = '  +0.12%  '

$ws.Range("D21").Value = '6.917'
$ws.Range("E21").Value = '  -5.59%  '

$ws.Range("D23").Value = '13.62'
$ws.Range("E23").Value = '  -4.95%  '

$ws.Range("D24").Value = '24.080.76'
$ws.Range("E24").Value = '  -3.84%  '

$ws.Range("D25").Value = '2.331'
$ws.Range("E25").Value = '  -2.80%  '

$ws.Range("D26").Value = '2.707'
$ws.Range("E26").Value = '  -8.13%  '

$ws.Range("D27").Value = '22.45'
$ws.Range("E27").Value = '  -4.79%  '

$ws.Range("D28").Value = '157.91'
$ws.Range("E28").Value = '  -3.09%  '

$ws.Range("E29").Value = '  +3.84%  '

$ws.Range("D30").Value = '141.36'
$ws.Range("E30").Value = '  -7.13%  '

$ws.Range("D31").Value = '5.276'
$ws.Range("E31").Value = '  -13.21%  '

$ws.Range("D32").Value = '2.457'
$ws.Range("E32").Value = '  -9.18%  '

$ws.Range("D33").Value = '1.819.71'
$ws.Range("E33").Value = '  -3.20%  '

$ws.Range("D34").Value = '6.964'
$ws.Range("E34").Value = '  -3.39%  '

$ws.Range("D35").Value = '0.08014'
$ws.Range("E35").Value = '  -7.33%  '

$ws.Range("E36").Value = '  -8.34%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.9587'
$ws.Range("E37").Value = '  -8.60%  '

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = '0.2698'
$ws.Range("E38").Value = '  -7.96%  '

$ws.Range("D39").Value = '0.09199'
$ws.Range("E39").Value = '  -5.33%  '

$ws.Range("D40").Value = '1.459'
$ws.Range("E40").Value = '  -1.42%  '

$ws.Range("D41").Value = '9.950'
$ws.Range("E41").Value = '  -9.70%  '

$ws.Range("D42").Value = '0.7602'
$ws.Range("E42").Value = '  -8.45%  '

$ws.Range("D43").Value = '13.08'
$ws.Range("E43").Value = '  -7.00%  '

$ws.Range("D44").Value = '16.00'
$ws.Range("E44").Value = '  -6.26%  '

$ws.Range("D45").Value = '0.6906'
$ws.Range("E45").Value = '  -6.50%  '

$ws.Range("D46").Value = '2.480'
$ws.Range("E46").Value = '  -7.96%  '

$ws.Range("D47").Value = '4.098'
$ws.Range("E47").Value = '  -3.58%  '

$ws.Range("D48").Value = '0.9999'
$ws.Range("E48").Value = '  +0.00%  '

$ws.Range("D49").Value = '0.08344'
$ws.Range("E49").Value = '  -9.56%  '

$ws.Range("B50").Value = 'Flow'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D50").Value = '1.263'
$ws.Range("E50").Value = '  -10.37%  '

$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = '133.73'
$ws.Range("E51").Value = '  -4.56%  '
